# test/controllers/i18n.xlsx edit
# - Row 2, column B ("key" for the 2nd data row) changes from "test" to "r2space"
# - The active selection moves from C3 to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B2 cell value (was "test", becomes "r2space").
$ws.Range("B2").Value = "r2space"

# Move/replace the sheet's selection to B7.
$ws.Range("B7").Select()
